# Listas sem duplicação de professores
# Move the 'MCT-3A-CAM' marker out of rows/columns where it duplicates
# with another entry, and place it into the correct day column instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "-"
$ws.Range("E2").Value = "[-, -, -, 'MCT-3A-CAM']"

$ws.Range("E3").Value = "[-, -, -, 'MCT-3A-CAM']"

$ws.Range("B6").Value = "['MCT-3A-CAM', -, -, -]"
$ws.Range("F6").Value = "-"

$ws.Range("B7").Value = "['MCT-3A-CAM', -, -, -]"
$ws.Range("F7").Value = "-"

$ws.Range("B8").Value = "-"
